# Update the line-loading result table ("pl_mw" / 380 kV case) in place.
# Columns A,F,I,K,L,M are untouched (index / all-zero columns); only the
# computed power-flow columns B,C,D,E,G,H,J,N,O for data rows 2-25 change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value2 = 0.9840413895909137
$ws.Range("C2").Value2 = 0.1851151416863388
$ws.Range("D2").Value2 = 0.6476995259981493
$ws.Range("E2").Value2 = 0.2645308497797458
$ws.Range("G2").Value2 = 0.6317406077371004
$ws.Range("H2").Value2 = 0.7303384626437293
$ws.Range("J2").Value2 = 0.1379778847859612
$ws.Range("N2").Value2 = 1.041204385562807
$ws.Range("O2").Value2 = 2.711130602301807

# row 3
$ws.Range("B3").Value2 = 0.8853465973357402
$ws.Range("C3").Value2 = 0.161701311268132
$ws.Range("D3").Value2 = 0.6369589940424021
$ws.Range("E3").Value2 = 0.259267791642479
$ws.Range("G3").Value2 = 0.6272734791357948
$ws.Range("H3").Value2 = 0.7329613008094924
$ws.Range("J3").Value2 = 0.1344947395986935
$ws.Range("N3").Value2 = 1.047568016703494
$ws.Range("O3").Value2 = 2.706808290012731

# row 4
$ws.Range("B4").Value2 = 0.8248152129607433
$ws.Range("C4").Value2 = 0.1472780473741864
$ws.Range("D4").Value2 = 0.6307011226130896
$ws.Range("E4").Value2 = 0.2561800331177864
$ws.Range("G4").Value2 = 0.6250341534078387
$ws.Range("H4").Value2 = 0.7349699749693883
$ws.Range("J4").Value2 = 0.132435385498674
$ws.Range("N4").Value2 = 1.051920591855264
$ws.Range("O4").Value2 = 2.70604315630456

# row 5
$ws.Range("B5").Value2 = 0.8001665124358794
$ws.Range("C5").Value2 = 0.141388953586727
$ws.Range("D5").Value2 = 0.6282357506068479
$ws.Range("E5").Value2 = 0.2549579032611433
$ws.Range("G5").Value2 = 0.624247903127582
$ws.Range("H5").Value2 = 0.7358886242278686
$ws.Range("J5").Value2 = 0.1316161076831079
$ws.Range("N5").Value2 = 1.053806402634244
$ws.Range("O5").Value2 = 2.706205535235455

# row 6
$ws.Range("B6").Value2 = 0.7960747607763778
$ws.Range("C6").Value2 = 0.1404103897544644
$ws.Range("D6").Value2 = 0.6278314985173381
$ws.Range("E6").Value2 = 0.2547571533918642
$ws.Range("G6").Value2 = 0.6241249653557901
$ws.Range("H6").Value2 = 0.7360472096955561
$ws.Range("J6").Value2 = 0.1314812701813324
$ws.Range("N6").Value2 = 1.054126314689405
$ws.Range("O6").Value2 = 2.70626111471401

# row 7
$ws.Range("B7").Value2 = 0.8244827159861643
$ws.Range("C7").Value2 = 0.1471986710916156
$ws.Range("D7").Value2 = 0.6306675304730334
$ws.Range("E7").Value2 = 0.25616340464601
$ws.Range("G7").Value2 = 0.6250230388290703
$ws.Range("H7").Value2 = 0.7349819589594091
$ws.Range("J7").Value2 = 0.1324242557757671
$ws.Range("N7").Value2 = 1.051945570489657
$ws.Range("O7").Value2 = 2.706043427310448

# row 8
$ws.Range("B8").Value2 = 0.9499981553958037
$ws.Range("C8").Value2 = 0.1770520213069915
$ws.Range("D8").Value2 = 0.6439263037140677
$ws.Range("E8").Value2 = 0.2626863045160377
$ws.Range("G8").Value2 = 0.6300956434648697
$ws.Range("H8").Value2 = 0.7311601352213302
$ws.Range("J8").Value2 = 0.1367604231630608
$ws.Range("N8").Value2 = 1.043306224559032
$ws.Range("O8").Value2 = 2.709247776745229

# row 9
$ws.Range("B9").Value2 = 1.196624128355666
$ws.Range("C9").Value2 = 0.2352085939646713
$ws.Range("D9").Value2 = 0.6725992910409104
$ws.Range("E9").Value2 = 0.2766194619764804
$ws.Range("G9").Value2 = 0.6440540354073931
$ws.Range("H9").Value2 = 0.7268280014248347
$ws.Range("J9").Value2 = 0.1458943567208877
$ws.Range("N9").Value2 = 1.029892377833406
$ws.Range("O9").Value2 = 2.730557678482654

# row 10
$ws.Range("B10").Value2 = 1.378075822618371
$ws.Range("C10").Value2 = 0.277688733989379
$ws.Range("D10").Value2 = 0.6952975235843439
$ws.Range("E10").Value2 = 0.2875549906703085
$ws.Range("G10").Value2 = 0.656778749274622
$ws.Range("H10").Value2 = 0.7255776386539736
$ws.Range("J10").Value2 = 0.1529924359873718
$ws.Range("N10").Value2 = 1.022181761578864
$ws.Range("O10").Value2 = 2.755435886447827

# row 11
$ws.Range("B11").Value2 = 1.46067035254913
$ws.Range("C11").Value2 = 0.2969579584543283
$ws.Range("D11").Value2 = 0.7059787932811616
$ws.Range("E11").Value2 = 0.2926822866923899
$ws.Range("G11").Value2 = 0.6631091314588247
$ws.Range("H11").Value2 = 0.7254294896419822
$ws.Range("J11").Value2 = 0.1563063265813724
$ws.Range("N11").Value2 = 1.019138527524134
$ws.Range("O11").Value2 = 2.768769743580691

# row 12
$ws.Range("B12").Value2 = 1.491953037554936
$ws.Range("C12").Value2 = 0.3042464851607463
$ws.Range("D12").Value2 = 0.7100746672278149
$ws.Range("E12").Value2 = 0.2946458392967841
$ws.Range("G12").Value2 = 0.6655845876675954
$ws.Range("H12").Value2 = 0.7254339508667016
$ws.Range("J12").Value2 = 0.1575734603589325
$ws.Range("N12").Value2 = 1.018052820528084
$ws.Range("O12").Value2 = 2.774109902631892

# row 13
$ws.Range("B13").Value2 = 1.485215510790852
$ws.Range("C13").Value2 = 0.3026771460421287
$ws.Range("D13").Value2 = 0.7091902751100463
$ws.Range("E13").Value2 = 0.294221976519438
$ws.Range("G13").Value2 = 0.6650479667048188
$ws.Range("H13").Value2 = 0.7254302953708134
$ws.Range("J13").Value2 = 0.1573000159188211
$ws.Range("N13").Value2 = 1.018283681604444
$ws.Range("O13").Value2 = 2.772946850459022

# row 14
$ws.Range("B14").Value2 = 1.463243887273222
$ws.Range("C14").Value2 = 0.2975577586113332
$ws.Range("D14").Value2 = 0.7063147392511269
$ws.Range("E14").Value2 = 0.2928433892789997
$ws.Range("G14").Value2 = 0.6633112178586771
$ws.Range("H14").Value2 = 0.7254286424928296
$ws.Range("J14").Value2 = 0.1564103291035224
$ws.Range("N14").Value2 = 1.019047869405568
$ws.Range("O14").Value2 = 2.769203245962132

# row 15
$ws.Range("B15").Value2 = 1.449786362948657
$ws.Range("C15").Value2 = 0.2944208927149532
$ws.Range("D15").Value2 = 0.7045600450827578
$ws.Range("E15").Value2 = 0.2920018243217797
$ws.Range("G15").Value2 = 0.6622576138436642
$ws.Range("H15").Value2 = 0.7254355192147983
$ws.Range("J15").Value2 = 0.1558669642736419
$ws.Range("N15").Value2 = 1.019524640633591
$ws.Range("O15").Value2 = 2.766948093587871

# row 16
$ws.Range("B16").Value2 = 1.372678983806509
$ws.Range("C16").Value2 = 0.2764283021609515
$ws.Range("D16").Value2 = 0.6946066318600685
$ws.Range("E16").Value2 = 0.2872229821831453
$ws.Range("G16").Value2 = 0.6563759781081586
$ws.Range("H16").Value2 = 0.7255957913040589
$ws.Range("J16").Value2 = 0.1527775764418777
$ws.Range("N16").Value2 = 1.022389981701906
$ws.Range("O16").Value2 = 2.754605146609634

# row 17
$ws.Range("B17").Value2 = 1.325388202935528
$ws.Range("C17").Value2 = 0.2653760262686546
$ws.Range("D17").Value2 = 0.6885916142943529
$ws.Range("E17").Value2 = 0.2843304217932783
$ws.Range("G17").Value2 = 0.6529068218315217
$ws.Range("H17").Value2 = 0.725801902964875
$ws.Range("J17").Value2 = 0.1509041070007555
$ws.Range("N17").Value2 = 1.024266654960201
$ws.Range("O17").Value2 = 2.747550316472569

# row 18
$ws.Range("B18").Value2 = 1.298192705631322
$ws.Range("C18").Value2 = 0.2590138755577698
$ws.Range("D18").Value2 = 0.6851654208978744
$ws.Range("E18").Value2 = 0.2826810673043028
$ws.Range("G18").Value2 = 0.6509624344328415
$ws.Range("H18").Value2 = 0.7259600418090599
$ws.Range("J18").Value2 = 0.149834530491475
$ws.Range("N18").Value2 = 1.025389780030451
$ws.Range("O18").Value2 = 2.743682324549809

# row 19
$ws.Range("B19").Value2 = 1.288985664624875
$ws.Range("C19").Value2 = 0.2568588833760543
$ws.Range("D19").Value2 = 0.6840111232380082
$ws.Range("E19").Value2 = 0.2821250921353311
$ws.Range("G19").Value2 = 0.6503128430510969
$ws.Range("H19").Value2 = 0.7260203821071229
$ws.Range("J19").Value2 = 0.1494737625189089
$ws.Range("N19").Value2 = 1.025777560895818
$ws.Range("O19").Value2 = 2.742405252153077

# row 20
$ws.Range("B20").Value2 = 1.330421889666241
$ws.Range("C20").Value2 = 0.2665530980887354
$ws.Range("D20").Value2 = 0.6892284578724173
$ws.Range("E20").Value2 = 0.2846368525484309
$ws.Range("G20").Value2 = 0.6532708406310661
$ws.Range("H20").Value2 = 0.7257758641947447
$ws.Range("J20").Value2 = 0.151102713564498
$ws.Range("N20").Value2 = 1.024062356395795
$ws.Range("O20").Value2 = 2.74828166903859

# row 21
$ws.Range("B21").Value2 = 1.469697336235754
$ws.Range("C21").Value2 = 0.2990616752953486
$ws.Range("D21").Value2 = 0.707157967108003
$ws.Range("E21").Value2 = 0.2932477178388382
$ws.Range("G21").Value2 = 0.6638192157787302
$ws.Range("H21").Value2 = 0.7254274837577981
$ws.Range("J21").Value2 = 0.1566713194542473
$ws.Range("N21").Value2 = 1.018821599284237
$ws.Range("O21").Value2 = 2.770294930559629

# row 22
$ws.Range("B22").Value2 = 1.560755786498532
$ws.Range("C22").Value2 = 0.3202592987539958
$ws.Range("D22").Value2 = 0.7191738194839559
$ws.Range("E22").Value2 = 0.2990033992236505
$ws.Range("G22").Value2 = 0.6711696194352186
$ws.Range("H22").Value2 = 0.7255528335138877
$ws.Range("J22").Value2 = 0.1603820521234098
$ws.Range("N22").Value2 = 1.015785207617839
$ws.Range("O22").Value2 = 2.786377985788505

# row 23
$ws.Range("B23").Value2 = 1.512153550850485
$ws.Range("C23").Value2 = 0.3089502986415198
$ws.Range("D23").Value2 = 0.7127334913371044
$ws.Range("E23").Value2 = 0.2959197721550169
$ws.Range("G23").Value2 = 0.6672046879812683
$ws.Range("H23").Value2 = 0.7254536046734898
$ws.Range("J23").Value2 = 0.158395030804158
$ws.Range("N23").Value2 = 1.017370240064551
$ws.Range("O23").Value2 = 2.777638655245198

# row 24
$ws.Range("B24").Value2 = 1.32814618282714
$ws.Range("C24").Value2 = 0.2660209689674957
$ws.Range("D24").Value2 = 0.688940441461682
$ws.Range("E24").Value2 = 0.2844982727902163
$ws.Range("G24").Value2 = 0.6531061117715211
$ws.Range("H24").Value2 = 0.7257875128405686
$ws.Range("J24").Value2 = 0.1510129001594578
$ws.Range("N24").Value2 = 1.024154582085103
$ws.Range("O24").Value2 = 2.747950439230436

# row 25
$ws.Range("B25").Value2 = 1.129857263837664
$ws.Range("C25").Value2 = 0.2195182778781088
$ws.Range("D25").Value2 = 0.6645560716613375
$ws.Range("E25").Value2 = 0.2727276524006541
$ws.Range("G25").Value2 = 0.6398460052765955
$ws.Range("H25").Value2 = 0.7276609297926768
$ws.Range("J25").Value2 = 0.1433555659683137
$ws.Range("N25").Value2 = 1.033144154222732
$ws.Range("O25").Value2 = 2.772946850459022
